$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in column B for all data rows (2-38): 18:41:30 -> 19:24:51
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 2).Value = "2023-06-26 19:24:51"
}

# Specific numeric corrections (columns E = Review Count, F = Rank)
$ws.Cells.Item(5, 6).Value = 117      # F5: 126 -> 117
$ws.Cells.Item(8, 5).Value = 320687   # E8: 322237 -> 320687
$ws.Cells.Item(10, 5).Value = 42773   # E10: 42796 -> 42773
$ws.Cells.Item(10, 6).Value = 84      # F10: 94 -> 84
$ws.Cells.Item(12, 5).Value = 2003619 # E12: 2010064 -> 2003619
$ws.Cells.Item(12, 6).Value = 13      # F12: 12 -> 13
$ws.Cells.Item(14, 6).Value = 111     # F14: 118 -> 111
$ws.Cells.Item(15, 6).Value = 196     # F15: (empty) -> 196
$ws.Cells.Item(22, 6).Value = 54      # F22: 58 -> 54
